$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Agrp"
$ws.Range("C2").Value = "Mc4r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2858313333333333
$ws.Range("H2").Value = 0.857494
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01244633333333333
$ws.Range("N2").Value = 0.037339
$ws.Range("O2").Value = 0.03341097222806328
$ws.Range("P2").Value = 0.03341097222806328
$ws.Range("Q2").Value = 0.003557552051777777
$ws.Range("R2").Value = 0.032017968466
$ws.Range("S2").Value = 0.03341097222806328
$ws.Range("T2").Value = 0.03341097222806328

$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Agrp"
$ws.Range("C3").Value = "Mc4r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2858313333333333
$ws.Range("H3").Value = 0.857494
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.05317933333333334
$ws.Range("N3").Value = 0.159538
$ws.Range("O3").Value = 0.1427547520640821
$ws.Range("P3").Value = 0.1427547520640821
$ws.Range("Q3").Value = 0.01520031975244444
$ws.Range("R3").Value = 0.136802877772
$ws.Range("S3").Value = 0.1427547520640821
$ws.Range("T3").Value = 0.1427547520640821

$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Agrp"
$ws.Range("C4").Value = "Mc4r"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2858313333333333
$ws.Range("H4").Value = 0.857494
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.282443
$ws.Range("N4").Value = 0.847329
$ws.Range("O4").Value = 0.758190784087218
$ws.Range("P4").Value = 0.7581907840872181
$ws.Range("Q4").Value = 0.08073105928066666
$ws.Range("R4").Value = 0.726579533526
$ws.Range("S4").Value = 0.758190784087218
$ws.Range("T4").Value = 0.7581907840872181

$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Agrp"
$ws.Range("C5").Value = "Mc4r"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2858313333333333
$ws.Range("H5").Value = 0.857494
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.02445366666666667
$ws.Range("N5").Value = 0.073361
$ws.Range("O5").Value = 0.06564349162063661
$ws.Range("P5").Value = 0.06564349162063661
$ws.Range("Q5").Value = 0.006989624148222222
$ws.Range("R5").Value = 0.06290661733399999
$ws.Range("S5").Value = 0.06564349162063661
$ws.Range("T5").Value = 0.06564349162063661
